$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Two stock codes changed name (stock-fishing country correction):
#   CODCOASTNOR -> CODNEARNCW
#   CODNS       -> CODIIIaW
# Update the corresponding STOCK cells in column B so the data (and the
# chart that reads from it) reflects the corrected stock codes.

$ws.Range("B11").Value = "CODNEARNCW"
$ws.Range("B11").Style = "Normal"

$ws.Range("B6").Value = "CODIIIaW"

[void]$ws.Range("B6").Select()
